$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 351 ("Femacal de La Calera" /
# Zapallo italiano weekly update). Everything that was on rows 351-356
# shifts down to 353-358, matching the diff.
$ws.Rows("351:352").Insert()

# New row 351
$ws.Range("A351").Value = 3
$ws.Range("B351").Value = "Femacal de La Calera"
$ws.Range("C351").Value = "Coquimbo"
$ws.Range("D351").Value = 44595
$ws.Range("E351").Value = 5
$ws.Range("F351").Value = 100112032
$ws.Range("G351").Value = "Zapallo italiano"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 160
$ws.Range("K351").Value = 4500
$ws.Range("L351").Value = 5000
$ws.Range("M351").Value = 4750
$ws.Range("N351").Value = "$/caja 36 unidades"
$ws.Range("O351").Value = "Provincia de Quillota"
$ws.Range("P351").Value = 132
$ws.Range("Q351").Value = 36
$ws.Range("R351").Value = "Hortaliza"

# New row 352
$ws.Range("A352").Value = 3
$ws.Range("B352").Value = "Femacal de La Calera"
$ws.Range("C352").Value = "Coquimbo"
$ws.Range("D352").Value = 44595
$ws.Range("E352").Value = 5
$ws.Range("F352").Value = 100112032
$ws.Range("G352").Value = "Zapallo italiano"
$ws.Range("H352").Value = "Sin especificar"
$ws.Range("I352").Value = "Primera"
$ws.Range("J352").Value = 160
$ws.Range("K352").Value = 9000
$ws.Range("L352").Value = 9500
$ws.Range("M352").Value = 9250
$ws.Range("N352").Value = "$/caja 70 unidades"
$ws.Range("O352").Value = "Región de Arica y Parinacota"
$ws.Range("P352").Value = 132
$ws.Range("Q352").Value = 70
$ws.Range("R352").Value = "Hortaliza"
